$wb = $excel.ActiveWorkbook

# The new sheet's data/layout/conditional-formatting are a near-identical copy
# of "GRASP (rec_seed)"; copying it preserves styles, conditional formatting,
# column widths, page setup, etc. exactly, then we overwrite the data values
# and metadata that differ.
$src = $wb.Worksheets.Item("GRASP (rec_seed)")
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "GRASP (randnode_rec_seed)"

$new.Cells.Item(2,1).Value2 = 0.1
$new.Cells.Item(2,2).Value2 = 20
$new.Cells.Item(2,3).Value2 = 42783759
$new.Cells.Item(2,4).Value2 = 133241309
$new.Cells.Item(2,5).Value2 = 7686
$new.Cells.Item(2,6).Value2 = 39144
$new.Cells.Item(2,7).Value2 = 22068
$new.Cells.Item(2,8).Value2 = 103247
$new.Cells.Item(2,9).Value2 = 95859
$new.Cells.Item(2,10).Value2 = 34316
$new.Cells.Item(2,11).Value2 = 80775
$new.Cells.Item(3,1).Value2 = 0.2
$new.Cells.Item(3,2).Value2 = 20
$new.Cells.Item(3,3).Value2 = 42630715
$new.Cells.Item(3,4).Value2 = 132987802
$new.Cells.Item(3,5).Value2 = 7713
$new.Cells.Item(3,6).Value2 = 40124
$new.Cells.Item(3,7).Value2 = 22068
$new.Cells.Item(3,8).Value2 = 106706
$new.Cells.Item(3,9).Value2 = 94829
$new.Cells.Item(3,10).Value2 = 34316
$new.Cells.Item(3,11).Value2 = 80775
$new.Cells.Item(4,1).Value2 = 0.3
$new.Cells.Item(4,2).Value2 = 20
$new.Cells.Item(4,3).Value2 = 42630715
$new.Cells.Item(4,4).Value2 = 132894610
$new.Cells.Item(4,5).Value2 = 7375
$new.Cells.Item(4,6).Value2 = 39274
$new.Cells.Item(4,7).Value2 = 22068
$new.Cells.Item(4,8).Value2 = 103674
$new.Cells.Item(4,9).Value2 = 95797
$new.Cells.Item(4,10).Value2 = 34316
$new.Cells.Item(4,11).Value2 = 80775
$new.Cells.Item(5,1).Value2 = 0.4
$new.Cells.Item(5,2).Value2 = 20
$new.Cells.Item(5,3).Value2 = 42085075
$new.Cells.Item(5,4).Value2 = 132894610
$new.Cells.Item(5,5).Value2 = 7511
$new.Cells.Item(5,6).Value2 = 39274
$new.Cells.Item(5,7).Value2 = 22068
$new.Cells.Item(5,8).Value2 = 106562
$new.Cells.Item(5,9).Value2 = 98072
$new.Cells.Item(5,10).Value2 = 34316
$new.Cells.Item(5,11).Value2 = 80775
$new.Cells.Item(6,1).Value2 = 0.5
$new.Cells.Item(6,2).Value2 = 20
$new.Cells.Item(6,3).Value2 = 42085075
$new.Cells.Item(6,4).Value2 = 132330304
$new.Cells.Item(6,5).Value2 = 7382
$new.Cells.Item(6,6).Value2 = 38304
$new.Cells.Item(6,7).Value2 = 22068
$new.Cells.Item(6,8).Value2 = 106151
$new.Cells.Item(6,9).Value2 = 102485
$new.Cells.Item(6,10).Value2 = 34316
$new.Cells.Item(6,11).Value2 = 80775
$new.Cells.Item(7,1).Value2 = 0.6
$new.Cells.Item(7,2).Value2 = 20
$new.Cells.Item(7,3).Value2 = 42630715
$new.Cells.Item(7,4).Value2 = 131022310
$new.Cells.Item(7,5).Value2 = 7382
$new.Cells.Item(7,6).Value2 = 36628
$new.Cells.Item(7,7).Value2 = 22068
$new.Cells.Item(7,8).Value2 = 106937
$new.Cells.Item(7,9).Value2 = 101839
$new.Cells.Item(7,10).Value2 = 34316
$new.Cells.Item(7,11).Value2 = 80343
$new.Cells.Item(8,1).Value2 = 0.7
$new.Cells.Item(8,2).Value2 = 20
$new.Cells.Item(8,3).Value2 = 42630715
$new.Cells.Item(8,4).Value2 = 130736205
$new.Cells.Item(8,5).Value2 = 7370
$new.Cells.Item(8,6).Value2 = 39214
$new.Cells.Item(8,7).Value2 = 22068
$new.Cells.Item(8,8).Value2 = 104879
$new.Cells.Item(8,9).Value2 = 102892
$new.Cells.Item(8,10).Value2 = 34316
$new.Cells.Item(8,11).Value2 = 80775
$new.Cells.Item(9,1).Value2 = 0.8
$new.Cells.Item(9,2).Value2 = 20
$new.Cells.Item(9,3).Value2 = 42661819
$new.Cells.Item(9,4).Value2 = 131719633
$new.Cells.Item(9,5).Value2 = 7425
$new.Cells.Item(9,6).Value2 = 40264
$new.Cells.Item(9,7).Value2 = 22068
$new.Cells.Item(9,8).Value2 = 108186
$new.Cells.Item(9,9).Value2 = 102892
$new.Cells.Item(9,10).Value2 = 34316
$new.Cells.Item(9,11).Value2 = 80775
$new.Cells.Item(10,1).Value2 = 0.9
$new.Cells.Item(10,2).Value2 = 20
$new.Cells.Item(10,3).Value2 = 42630715
$new.Cells.Item(10,4).Value2 = 131022310
$new.Cells.Item(10,5).Value2 = 7586
$new.Cells.Item(10,6).Value2 = 39314
$new.Cells.Item(10,7).Value2 = 22068
$new.Cells.Item(10,8).Value2 = 109672
$new.Cells.Item(10,9).Value2 = 102892
$new.Cells.Item(10,10).Value2 = 34316
$new.Cells.Item(10,11).Value2 = 80775
$new.Cells.Item(11,1).Value2 = 1
$new.Cells.Item(11,2).Value2 = 20
$new.Cells.Item(11,3).Value2 = 42783759
$new.Cells.Item(11,4).Value2 = 131022310
$new.Cells.Item(11,5).Value2 = 7568
$new.Cells.Item(11,6).Value2 = 42164
$new.Cells.Item(11,7).Value2 = 22068
$new.Cells.Item(11,8).Value2 = 102892
$new.Cells.Item(11,9).Value2 = 102892
$new.Cells.Item(11,10).Value2 = 34316
$new.Cells.Item(11,11).Value2 = 80775

# Update the selection on the source sheet (it's no longer the active tab)
# and on the new sheet (which becomes the active tab).
$src.Range("G59").Select()
$new.Range("I11").Select()
